$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Volume) updates: plain text assignments ---
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E33").Value = "  -4.76%  "
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  +10.38%  "
$ws.Range("E51").Value = "  -1.22%  "

# --- Column D (Price) updates that remain non-numeric text: plain text assignments ---
$ws.Range("D2").Value = "60.979.34"
$ws.Range("D3").Value = "3.383.73"
$ws.Range("D12").Value = "3.961.95"
$ws.Range("D16").Value = "3.387.95"
$ws.Range("D17").Value = "61.069.46"
$ws.Range("D26").Value = "3.522.72"
$ws.Range("D37").Value = "3.415.90"
$ws.Range("D47").Value = "2.444.42"

# --- Column D (Price) updates that look like plain numbers: must force text ---
# Using a text-formula + copy/paste-special-values trick avoids Excel converting
# the numeric-looking string into a real number (which would drop the formatting,
# e.g. turn "1.00" into 1) and also avoids creating new cell styles.
$ws.Range("D4").Formula = "=""0.999"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("D5").Formula = "=""571.04"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("D6").Formula = "=""141.95"""
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("D14").Formula = "=""27.84"""
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("D18").Formula = "=""6.08"""
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("D22").Formula = "=""75.24"""
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("D24").Formula = "=""1.00"""
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("D28").Formula = "=""0.999"""
$ws.Range("D28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163) | Out-Null
$ws.Range("D29").Formula = "=""7.22"""
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("D33").Formula = "=""1.38"""
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("D34").Formula = "=""23.19"""
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("D36").Formula = "=""166.28"""
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("D38").Formula = "=""4.98"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("D41").Formula = "=""26.84"""
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("D48").Formula = "=""22.90"""
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false
